$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the edited range to Text format so numeric-looking strings
# (e.g. "520.86", "1.00") are stored as literal text, matching the source data.
$editRange = $ws.Range("D2:E51")
$editRange.NumberFormat = "@"

$ws.Range('D2').Value = '58.059.46'
$ws.Range('E2').Value = '  -1.35%  '
$ws.Range('D3').Value = '2.464.60'
$ws.Range('E3').Value = '  -1.21%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '520.86'
$ws.Range('E5').Value = '  -2.37%  '
$ws.Range('D6').Value = '133.26'
$ws.Range('E6').Value = '  -2.04%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '0.557'
$ws.Range('E8').Value = '  -1.38%  '
$ws.Range('D9').Value = '2.472.59'
$ws.Range('E9').Value = '  -1.35%  '
$ws.Range('D10').Value = '0.0979'
$ws.Range('E10').Value = '  -3.11%  '
$ws.Range('D11').Value = '0.157'
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('D12').Value = '5.29'
$ws.Range('E12').Value = '  -2.00%  '
$ws.Range('D13').Value = '0.339'
$ws.Range('E13').Value = '  -2.58%  '
$ws.Range('D14').Value = '2.906.52'
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('D15').Value = '57.984.49'
$ws.Range('E15').Value = '  -1.24%  '
$ws.Range('D16').Value = '22.29'
$ws.Range('E16').Value = '  -2.33%  '
$ws.Range('E17').Value = '  -2.16%  '
$ws.Range('D18').Value = '2.467.20'
$ws.Range('E18').Value = '  -1.44%  '
$ws.Range('D19').Value = '10.63'
$ws.Range('E19').Value = '  -3.62%  '
$ws.Range('D20').Value = '320.67'
$ws.Range('E20').Value = '  -0.57%  '
$ws.Range('D21').Value = '4.16'
$ws.Range('E21').Value = '  -1.93%  '
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = '5.72'
$ws.Range('E23').Value = '  -4.03%  '
$ws.Range('D24').Value = '64.73'
$ws.Range('E24').Value = '  -0.78%  '
$ws.Range('D25').Value = '0.409'
$ws.Range('E25').Value = '  -2.62%  '
$ws.Range('E27').Value = '  -2.87%  '
$ws.Range('D28').Value = '7.34'
$ws.Range('E28').Value = '  -2.15%  '
$ws.Range('D29').Value = '0.0₃0747'
$ws.Range('E29').Value = '  -1.90%  '
$ws.Range('D30').Value = '167.58'
$ws.Range('E30').Value = '  -2.47%  '
$ws.Range('E31').Value = '  -2.88%  '
$ws.Range('D32').Value = '6.23'
$ws.Range('E32').Value = '  -4.42%  '
$ws.Range('E33').Value = '  +0.13%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').Value = '1.36'
$ws.Range('E36').Value = '  +0.41%  '
$ws.Range('D37').Value = '18.04'
$ws.Range('E37').Value = '  -1.41%  '
$ws.Range('D38').Value = '3.97'
$ws.Range('E38').Value = '  -1.65%  '
$ws.Range('D39').Value = '36.28'
$ws.Range('E39').Value = '  -1.30%  '
$ws.Range('D40').Value = '1.47'
$ws.Range('E40').Value = '  -4.07%  '
$ws.Range('D41').Value = '0.794'
$ws.Range('E41').Value = '  -1.23%  '
$ws.Range('D42').Value = '3.45'
$ws.Range('E42').Value = '  -3.22%  '
$ws.Range('D43').Value = '272.93'
$ws.Range('E43').Value = '  -3.30%  '
$ws.Range('E44').Value = '  -4.02%  '
$ws.Range('D45').Value = '0.589'
$ws.Range('E45').Value = '  -2.71%  '
$ws.Range('D46').Value = '124.52'
$ws.Range('E46').Value = '  -4.13%  '
$ws.Range('D47').Value = '0.0907'
$ws.Range('E47').Value = '  -1.61%  '
$ws.Range('D48').Value = '0.0487'
$ws.Range('E48').Value = '  -2.88%  '
$ws.Range('E49').Value = '  -2.71%  '
$ws.Range('D50').Value = '16.86'
$ws.Range('E50').Value = '  -2.10%  '
$ws.Range('D51').Value = '1.724.27'
$ws.Range('E51').Value = '  -1.69%  '

# Restore default (no explicit number format) so styling matches the original file.
$editRange.ClearFormats()
